# Update countries & provincias Spain
#
# Refreshes the COVID-19 "Pais" sheet with a newer data pull:
#   - bumps the "last updated" timestamp (06:46 -> 07:16)
#   - updates totals for Estados Unidos (row 6) and Austria (row 15)
#   - Armenia's case count overtakes Eslovaquia/Kuwait/Serbia/Bulgaria, so it
#     moves up to row 62, pushing those four down one row each (rows 62-66)
#   - Georgia's case count overtakes Camboya/Azerbaiyan/Estado de
#     Palestina/Oman/Trinidad yTobago, so it moves up to row 98, pushing
#     those five down one row each (rows 98-103)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "last updated" timestamp ---------------------------------------------
$ws.Cells.Item(1,1).Value = 'Datos actualizados a 22 de Marzo de 2020 a las 07:16'

# --- straightforward value refreshes (no re-ranking) -----------------------
# Row 6: Estados Unidos
$ws.Cells.Item(6,2).Value = 26888
$ws.Cells.Item(6,3).Value = 2681
$ws.Cells.Item(6,5).Value = 26362

# Row 15: Austria
$ws.Cells.Item(15,2).Value = 3024
$ws.Cells.Item(15,3).Value = 32
$ws.Cells.Item(15,5).Value = 3007

# --- Armenia re-rank block (rows 62-66) ------------------------------------
# Row 62: Armenia (new entry at this rank)
$ws.Cells.Item(62,1).Value = 'Armenia'
$ws.Cells.Item(62,2).Value = 190
$ws.Cells.Item(62,3).Value = 30
$ws.Cells.Item(62,4).Value = 2
$ws.Cells.Item(62,5).Value = 188
$ws.Cells.Item(62,6).Value = 2
$ws.Cells.Item(62,7).Value = 0
$ws.Cells.Item(62,8).Value = 0

# Row 63: Eslovaquia (shifted down from row 62)
$ws.Cells.Item(63,1).Value = 'Eslovaquia'
$ws.Cells.Item(63,2).Value = 178
$ws.Cells.Item(63,3).Value = 0
$ws.Cells.Item(63,4).Value = 7
$ws.Cells.Item(63,5).Value = 171
$ws.Cells.Item(63,6).Value = 2
$ws.Cells.Item(63,7).Value = 0
$ws.Cells.Item(63,8).Value = 0

# Row 64: Kuwait (shifted down from row 63)
$ws.Cells.Item(64,1).Value = 'Kuwait'
$ws.Cells.Item(64,2).Value = 176
$ws.Cells.Item(64,3).Value = 0
$ws.Cells.Item(64,4).Value = 27
$ws.Cells.Item(64,5).Value = 149
$ws.Cells.Item(64,6).Value = 5
$ws.Cells.Item(64,7).Value = 0
$ws.Cells.Item(64,8).Value = 0

# Row 65: Serbia (shifted down from row 64)
$ws.Cells.Item(65,1).Value = 'Serbia'
$ws.Cells.Item(65,2).Value = 171
$ws.Cells.Item(65,3).Value = 0
$ws.Cells.Item(65,4).Value = 2
$ws.Cells.Item(65,5).Value = 168
$ws.Cells.Item(65,6).Value = 4
$ws.Cells.Item(65,7).Value = 0
$ws.Cells.Item(65,8).Value = 1

# Row 66: Bulgaria (shifted down from row 65)
$ws.Cells.Item(66,1).Value = 'Bulgaria'
$ws.Cells.Item(66,2).Value = 163
$ws.Cells.Item(66,3).Value = 0
$ws.Cells.Item(66,4).Value = 3
$ws.Cells.Item(66,5).Value = 157
$ws.Cells.Item(66,6).Value = 3
$ws.Cells.Item(66,7).Value = 0
$ws.Cells.Item(66,8).Value = 3

# Row 67 (San Marino) is unaffected by the re-rank; left untouched.

# --- Georgia re-rank block (rows 98-103) -----------------------------------
# Row 98: Georgia (new entry at this rank)
$ws.Cells.Item(98,1).Value = 'Georgia'
$ws.Cells.Item(98,2).Value = 54
$ws.Cells.Item(98,3).Value = 5
$ws.Cells.Item(98,4).Value = 1
$ws.Cells.Item(98,5).Value = 53
$ws.Cells.Item(98,6).Value = 1
$ws.Cells.Item(98,7).Value = 0
$ws.Cells.Item(98,8).Value = 0

# Row 99: Camboya (shifted down from row 98)
$ws.Cells.Item(99,1).Value = 'Camboya'
$ws.Cells.Item(99,2).Value = 53
$ws.Cells.Item(99,3).Value = 0
$ws.Cells.Item(99,4).Value = 2
$ws.Cells.Item(99,5).Value = 51
$ws.Cells.Item(99,6).Value = 0
$ws.Cells.Item(99,7).Value = 0
$ws.Cells.Item(99,8).Value = 0

# Row 100: Azerbaiyan (shifted down from row 99)
$ws.Cells.Item(100,1).Value = 'Azerbaiyan'
$ws.Cells.Item(100,2).Value = 53
$ws.Cells.Item(100,3).Value = 0
$ws.Cells.Item(100,4).Value = 11
$ws.Cells.Item(100,5).Value = 41
$ws.Cells.Item(100,6).Value = 0
$ws.Cells.Item(100,7).Value = 0
$ws.Cells.Item(100,8).Value = 1

# Row 101: Estado de Palestina (shifted down from row 100)
$ws.Cells.Item(101,1).Value = 'Estado de Palestina'
$ws.Cells.Item(101,2).Value = 53
$ws.Cells.Item(101,3).Value = 0
$ws.Cells.Item(101,4).Value = 17
$ws.Cells.Item(101,5).Value = 36
$ws.Cells.Item(101,6).Value = 0
$ws.Cells.Item(101,7).Value = 0
$ws.Cells.Item(101,8).Value = 0

# Row 102: Oman (shifted down from row 101)
$ws.Cells.Item(102,1).Value = 'Oman'
$ws.Cells.Item(102,2).Value = 52
$ws.Cells.Item(102,3).Value = 0
$ws.Cells.Item(102,4).Value = 13
$ws.Cells.Item(102,5).Value = 39
$ws.Cells.Item(102,6).Value = 0
$ws.Cells.Item(102,7).Value = 0
$ws.Cells.Item(102,8).Value = 0

# Row 103: Trinidad yTobago (shifted down from row 102)
$ws.Cells.Item(103,1).Value = 'Trinidad yTobago'
$ws.Cells.Item(103,2).Value = 49
$ws.Cells.Item(103,3).Value = 0
$ws.Cells.Item(103,4).Value = 0
$ws.Cells.Item(103,5).Value = 49
$ws.Cells.Item(103,6).Value = 0
$ws.Cells.Item(103,7).Value = 0
$ws.Cells.Item(103,8).Value = 0

# Row 104 (Reunion) is unaffected by the re-rank; left untouched.
